$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "eduardopais3@gmail.com"
$ws.Range("B5").Value = "Eduardo1234"

$ws.Range("A6").Value = "dadasdasda@gmail.com"
$ws.Range("B6").Value = "sadasdA1"
